$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column C (NIP) before the existing bid column, shifting bid -> D.
# The engine carries the left-neighbour (bidder) style onto the freshly inserted
# blank cells, so explicitly clear formatting on the new column first.
$ws.Columns.Item(3).Insert()
$ws.Range("C1:C10").ClearFormats()

# Header for new column
$ws.Cells.Item(1,3).Value = 'nip'

# Fill NIP values for existing rows 2-10
$ws.Cells.Item(2,3).Value = 5252516964
$ws.Cells.Item(3,3).Value = 7720100638
$ws.Cells.Item(4,3).Value = 8982126450
$ws.Cells.Item(5,3).Value = 7720100638
$ws.Cells.Item(6,3).Value = 8982126450
$ws.Cells.Item(7,3).Value = 5252516964
$ws.Cells.Item(8,3).Value = 8982126450
$ws.Cells.Item(9,3).Value = 5252516964
$ws.Cells.Item(10,3).Value = 7720100638

# Insert new rows 11-28 (inherits formatting, e.g. bidder-name style, from row above)
$ws.Range("A11:D28").Insert()

# Populate the new rows
$ws.Cells.Item(11,1).Value = 4
$ws.Cells.Item(11,2).Value = 'JONESBORO SPÓŁKA Z OGRANICZONĄ ODPOWIEDZIALNOŚCIĄ'
$ws.Cells.Item(11,3).Value = 5252516964
$ws.Cells.Item(11,4).Value = 3800

$ws.Cells.Item(12,1).Value = 4
$ws.Cells.Item(12,2).Value = '"KODREWEX" SPÓŁKA Z OGRANICZONĄ ODPOWIEDZIALNOŚCIĄ'
$ws.Cells.Item(12,3).Value = 7720100638
$ws.Cells.Item(12,4).Value = 3900

$ws.Cells.Item(13,1).Value = 5
$ws.Cells.Item(13,2).Value = 'JONESBORO SPÓŁKA Z OGRANICZONĄ ODPOWIEDZIALNOŚCIĄ'
$ws.Cells.Item(13,3).Value = 5252516964
$ws.Cells.Item(13,4).Value = 3800

$ws.Cells.Item(14,1).Value = 5
$ws.Cells.Item(14,2).Value = '"KODREWEX" SPÓŁKA Z OGRANICZONĄ ODPOWIEDZIALNOŚCIĄ'
$ws.Cells.Item(14,3).Value = 7720100638
$ws.Cells.Item(14,4).Value = 3900

$ws.Cells.Item(15,1).Value = 6
$ws.Cells.Item(15,2).Value = 'JONESBORO SPÓŁKA Z OGRANICZONĄ ODPOWIEDZIALNOŚCIĄ'
$ws.Cells.Item(15,3).Value = 5252516964
$ws.Cells.Item(15,4).Value = 3800

$ws.Cells.Item(16,1).Value = 6
$ws.Cells.Item(16,2).Value = '"KODREWEX" SPÓŁKA Z OGRANICZONĄ ODPOWIEDZIALNOŚCIĄ'
$ws.Cells.Item(16,3).Value = 7720100638
$ws.Cells.Item(16,4).Value = 3900

$ws.Cells.Item(17,1).Value = 7
$ws.Cells.Item(17,2).Value = 'JONESBORO SPÓŁKA Z OGRANICZONĄ ODPOWIEDZIALNOŚCIĄ'
$ws.Cells.Item(17,3).Value = 5252516964
$ws.Cells.Item(17,4).Value = 3800

$ws.Cells.Item(18,1).Value = 7
$ws.Cells.Item(18,2).Value = '"KODREWEX" SPÓŁKA Z OGRANICZONĄ ODPOWIEDZIALNOŚCIĄ'
$ws.Cells.Item(18,3).Value = 7720100638
$ws.Cells.Item(18,4).Value = 3900

$ws.Cells.Item(19,1).Value = 8
$ws.Cells.Item(19,2).Value = 'JONESBORO SPÓŁKA Z OGRANICZONĄ ODPOWIEDZIALNOŚCIĄ'
$ws.Cells.Item(19,3).Value = 5252516964
$ws.Cells.Item(19,4).Value = 3800

$ws.Cells.Item(20,1).Value = 8
$ws.Cells.Item(20,2).Value = '"KODREWEX" SPÓŁKA Z OGRANICZONĄ ODPOWIEDZIALNOŚCIĄ'
$ws.Cells.Item(20,3).Value = 7720100638
$ws.Cells.Item(20,4).Value = 5100

$ws.Cells.Item(21,1).Value = 9
$ws.Cells.Item(21,2).Value = 'JONESBORO SPÓŁKA Z OGRANICZONĄ ODPOWIEDZIALNOŚCIĄ'
$ws.Cells.Item(21,3).Value = 5252516964
$ws.Cells.Item(21,4).Value = 5200

$ws.Cells.Item(22,1).Value = 9
$ws.Cells.Item(22,2).Value = '"KODREWEX" SPÓŁKA Z OGRANICZONĄ ODPOWIEDZIALNOŚCIĄ'
$ws.Cells.Item(22,3).Value = 7720100638
$ws.Cells.Item(22,4).Value = 5300

$ws.Cells.Item(23,1).Value = 10
$ws.Cells.Item(23,2).Value = '"KODREWEX" SPÓŁKA Z OGRANICZONĄ ODPOWIEDZIALNOŚCIĄ'
$ws.Cells.Item(23,3).Value = 7720100638
$ws.Cells.Item(23,4).Value = 5400

$ws.Cells.Item(24,1).Value = 11
$ws.Cells.Item(24,2).Value = '"KODREWEX" SPÓŁKA Z OGRANICZONĄ ODPOWIEDZIALNOŚCIĄ'
$ws.Cells.Item(24,3).Value = 7720100638
$ws.Cells.Item(24,4).Value = 5500

$ws.Cells.Item(25,1).Value = 12
$ws.Cells.Item(25,2).Value = '"KODREWEX" SPÓŁKA Z OGRANICZONĄ ODPOWIEDZIALNOŚCIĄ'
$ws.Cells.Item(25,3).Value = 7720100638
$ws.Cells.Item(25,4).Value = 5600

$ws.Cells.Item(26,1).Value = 13
$ws.Cells.Item(26,2).Value = '"KODREWEX" SPÓŁKA Z OGRANICZONĄ ODPOWIEDZIALNOŚCIĄ'
$ws.Cells.Item(26,3).Value = 7720100638
$ws.Cells.Item(26,4).Value = 5700

$ws.Cells.Item(27,1).Value = 14
$ws.Cells.Item(27,2).Value = '"KODREWEX" SPÓŁKA Z OGRANICZONĄ ODPOWIEDZIALNOŚCIĄ'
$ws.Cells.Item(27,3).Value = 7720100638
$ws.Cells.Item(27,4).Value = 5800

$ws.Cells.Item(28,1).Value = 15
$ws.Cells.Item(28,2).Value = '"KODREWEX" SPÓŁKA Z OGRANICZONĄ ODPOWIEDZIALNOŚCIĄ'
$ws.Cells.Item(28,3).Value = 7720100638
$ws.Cells.Item(28,4).Value = 5000

# Update selection to match the final state of the workbook
$ws.Range("D28").Select()
